$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: extend the thin header/rule row into the new column (empty cell, same style as Q2) ---
$ws.Range("Q2").Copy($ws.Range("R2"))

# --- Row 3: year header "2021" (same style as the other year cells) ---
$ws.Range("Q3").Copy($ws.Range("R3"))
$ws.Range("R3").Value2 = 2021

# --- Row 4: bold summary row, needs its own number-format-on-bold-font style ---
$ws.Range("Q4").Copy($ws.Range("R4"))
$ws.Range("R4").Value2 = 18
$ws.Range("R4").NumberFormat = "0.0"

# --- Rows 5-12: data rows, formatted like column O (numeric, vertical-centered) ---
$ws.Range("O5").Copy($ws.Range("R5"))
$ws.Range("R5").Value2 = 1.7480265877296817

$ws.Range("O6").Copy($ws.Range("R6"))
$ws.Range("R6").Value2 = 4.1112601249414027

$ws.Range("O7").Copy($ws.Range("R7"))
$ws.Range("R7").Value2 = 1.5225742120245318

$ws.Range("O8").Copy($ws.Range("R8"))
$ws.Range("R8").Value2 = 1.2326518235454269

$ws.Range("O9").Copy($ws.Range("R9"))
$ws.Range("R9").Value2 = 4.0865392096984241

$ws.Range("O10").Copy($ws.Range("R10"))
$ws.Range("R10").Value2 = 0.84876624403485645

$ws.Range("O11").Copy($ws.Range("R11"))
$ws.Range("R11").Value2 = 2.1456657699653627

$ws.Range("O12").Copy($ws.Range("R12"))
$ws.Range("R12").Value2 = 1.8214779402142154

# --- Row 13: bottom (thick-bordered) total row ---
$ws.Range("O13").Copy($ws.Range("R13"))
$ws.Range("R13").Value2 = 0.51989507542472779

# --- Selection, as left by the editor after the edit ---
$null = $ws.Range("R24:R25").Select()
